$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates: force text storage to match original inline-string cells
# (avoid Excel's automatic number inference for numeric-looking text), then restore
# the default (unstyled) cell format so no stray style index is introduced.
$dCells = @("D2", "D3", "D5", "D8", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D23", "D24", "D25", "D28", "D30", "D32", "D33", "D34", "D36", "D38", "D39", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "30.680.31"
$ws.Range("D3").Value = "1.677.83"
$ws.Range("D5").Value = "219.97"
$ws.Range("D8").Value = "29.59"
$ws.Range("D10").Value = "0.0641"
$ws.Range("D11").Value = "0.0906"
$ws.Range("D12").Value = "1.918.65"
$ws.Range("D13").Value = "1.672.63"
$ws.Range("D14").Value = "0.609"
$ws.Range("D15").Value = "10.01"
$ws.Range("D16").Value = "4.04"
$ws.Range("D17").Value = "30.688.21"
$ws.Range("D18").Value = "66.37"
$ws.Range("D19").Value = "242.87"
$ws.Range("D23").Value = "10.03"
$ws.Range("D24").Value = "2.15"
$ws.Range("D25").Value = "159.11"
$ws.Range("D28").Value = "6.70"
$ws.Range("D30").Value = "0.0494"
$ws.Range("D32").Value = "3.47"
$ws.Range("D33").Value = "3.32"
$ws.Range("D34").Value = "1.502.87"
$ws.Range("D36").Value = "84.25"
$ws.Range("D38").Value = "0.602"
$ws.Range("D39").Value = "0.0178"
$ws.Range("D47").Value = "5.56"
$ws.Range("D48").Value = "50.76"
$ws.Range("D49").Value = "1.810.42"
$ws.Range("D50").Value = "94.02"
$ws.Range("D51").Value = "0.0₆0116"

foreach ($addr in $dCells) { $ws.Range($addr).Style = "Normal" }

# Column E (Volume 1h) updates
$ws.Range("E2").Value = "  +2.54%  "
$ws.Range("E3").Value = "  +2.76%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +2.47%  "
$ws.Range("E6").Value = "  +2.54%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +3.13%  "
$ws.Range("E9").Value = "  +2.66%  "
$ws.Range("E10").Value = "  +5.36%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("E12").Value = "  +2.76%  "
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("E14").Value = "  +8.20%  "
$ws.Range("E15").Value = "  +7.88%  "
$ws.Range("E16").Value = "  +5.15%  "
$ws.Range("E17").Value = "  +2.49%  "
$ws.Range("E18").Value = "  +3.51%  "
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("E20").Value = "  +3.02%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  +2.97%  "
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("E26").Value = "  +3.55%  "
$ws.Range("E27").Value = "  +2.24%  "
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("E31").Value = "  +3.29%  "
$ws.Range("E32").Value = "  +2.63%  "
$ws.Range("E33").Value = "  +4.64%  "
$ws.Range("E34").Value = "  +5.14%  "
$ws.Range("E35").Value = "  +7.80%  "
$ws.Range("E36").Value = "  +11.29%  "
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("E38").Value = "  +9.18%  "
$ws.Range("E39").Value = "  +4.99%  "
$ws.Range("E40").Value = "  -3.00%  "
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("E44").Value = "  +1.94%  "
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("E47").Value = "  +3.82%  "
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("E49").Value = "  +2.04%  "
$ws.Range("E50").Value = "  +4.15%  "
$ws.Range("E51").Value = "  +1.84%  "
